$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # ALC
$ws2 = $wb.Worksheets.Item(2)  # ARM
$ws3 = $wb.Worksheets.Item(3)  # BSM
$ws4 = $wb.Worksheets.Item(4)  # CRP
$ws5 = $wb.Worksheets.Item(5)  # CUL
$ws6 = $wb.Worksheets.Item(6)  # GSM
$ws8 = $wb.Worksheets.Item(8)  # WVR

# ALC row 17
$ws1.Cells.Item(17, 8).Value = 3996
$ws1.Cells.Item(17, 9).Value = 0
$ws1.Cells.Item(17, 10).Value = 3996
$ws1.Cells.Item(17, 11).Value = 0
$ws1.Cells.Item(17, 12).Value = 11988
$ws1.Cells.Item(17, 13).ClearContents()
$ws1.Cells.Item(17, 14).Value = -12324

# ALC row 40
$ws1.Cells.Item(40, 8).Value = 254324.83
$ws1.Cells.Item(40, 9).Value = 2737.375
$ws1.Cells.Item(40, 10).Value = 757499.75
$ws1.Cells.Item(40, 11).Value = 2737.375
$ws1.Cells.Item(40, 12).Value = 757499.75
$ws1.Cells.Item(40, 13).Value = -2562.375
$ws1.Cells.Item(40, 14).Value = -757849.75

# ALC row 86
$ws1.Cells.Item(86, 8).Value = 4016.818
$ws1.Cells.Item(86, 9).Value = 4666.6665
$ws1.Cells.Item(86, 10).Value = 3773.125
$ws1.Cells.Item(86, 11).Value = 4666.6665
$ws1.Cells.Item(86, 12).Value = 3773.125
$ws1.Cells.Item(86, 13).Value = -3543.6665
$ws1.Cells.Item(86, 14).Value = -6019.125

# ALC row 89
$ws1.Cells.Item(89, 8).Value = 4016.818
$ws1.Cells.Item(89, 9).Value = 4666.6665
$ws1.Cells.Item(89, 10).Value = 3773.125
$ws1.Cells.Item(89, 11).Value = 23333.3325
$ws1.Cells.Item(89, 12).Value = 18865.625
$ws1.Cells.Item(89, 13).Value = -17717.3325
$ws1.Cells.Item(89, 14).Value = -30097.625

# ALC row 106
$ws1.Cells.Item(106, 8).Value = 3304.5715
$ws1.Cells.Item(106, 9).Value = 2688.6667
$ws1.Cells.Item(106, 10).Value = 7000
$ws1.Cells.Item(106, 11).Value = 2688.6667
$ws1.Cells.Item(106, 12).Value = 7000
$ws1.Cells.Item(106, 13).Value = -2057.6667
$ws1.Cells.Item(106, 14).Value = -8262

# ALC row 107
$ws1.Cells.Item(107, 8).Value = 880.25
$ws1.Cells.Item(107, 9).Value = 942.1818
$ws1.Cells.Item(107, 10).Value = 199
$ws1.Cells.Item(107, 11).Value = 942.1818
$ws1.Cells.Item(107, 12).Value = 199
$ws1.Cells.Item(107, 13).Value = 977.8182
$ws1.Cells.Item(107, 14).Value = -4039

# ARM row 32
$ws2.Cells.Item(32, 8).Value = 9479.860000000001
$ws2.Cells.Item(32, 9).Value = 9080.825000000001
$ws2.Cells.Item(32, 10).Value = 14800.333
$ws2.Cells.Item(32, 11).Value = 9080.825000000001
$ws2.Cells.Item(32, 12).Value = 14800.333
$ws2.Cells.Item(32, 13).Value = -8793.825000000001
$ws2.Cells.Item(32, 14).Value = -15374.333

# ARM row 61
$ws2.Cells.Item(61, 8).Value = 3047.6667
$ws2.Cells.Item(61, 9).Value = 3047.6667
$ws2.Cells.Item(61, 10).Value = 0
$ws2.Cells.Item(61, 11).Value = 3047.6667
$ws2.Cells.Item(61, 12).Value = 0
$ws2.Cells.Item(61, 13).Value = -2835.6667

# ARM row 63
$ws2.Cells.Item(63, 8).Value = 6314.857
$ws2.Cells.Item(63, 9).Value = 839.8
$ws2.Cells.Item(63, 10).Value = 20002.5
$ws2.Cells.Item(63, 11).Value = 839.8
$ws2.Cells.Item(63, 12).Value = 20002.5
$ws2.Cells.Item(63, 13).Value = -153.8
$ws2.Cells.Item(63, 14).Value = -21374.5

# ARM row 66
$ws2.Cells.Item(66, 8).Value = 6314.857
$ws2.Cells.Item(66, 9).Value = 839.8
$ws2.Cells.Item(66, 10).Value = 20002.5
$ws2.Cells.Item(66, 11).Value = 4199
$ws2.Cells.Item(66, 12).Value = 100012.5
$ws2.Cells.Item(66, 13).Value = -767
$ws2.Cells.Item(66, 14).Value = -106876.5

# ARM row 74
$ws2.Cells.Item(74, 8).Value = 8368.875
$ws2.Cells.Item(74, 9).Value = 8368.875
$ws2.Cells.Item(74, 10).Value = 0
$ws2.Cells.Item(74, 11).Value = 8368.875
$ws2.Cells.Item(74, 12).Value = 0
$ws2.Cells.Item(74, 13).Value = -7494.875

# ARM row 77
$ws2.Cells.Item(77, 8).Value = 8368.875
$ws2.Cells.Item(77, 9).Value = 8368.875
$ws2.Cells.Item(77, 10).Value = 0
$ws2.Cells.Item(77, 11).Value = 41844.375
$ws2.Cells.Item(77, 12).Value = 0
$ws2.Cells.Item(77, 13).Value = -37476.375

# ARM row 122
$ws2.Cells.Item(122, 8).Value = 1817.5
$ws2.Cells.Item(122, 9).Value = 1588.75
$ws2.Cells.Item(122, 10).Value = 3190
$ws2.Cells.Item(122, 11).Value = 4766.25
$ws2.Cells.Item(122, 12).Value = 9570
$ws2.Cells.Item(122, 13).Value = -2316.25
$ws2.Cells.Item(122, 14).Value = -14470

# ARM row 132
$ws2.Cells.Item(132, 8).Value = 4773
$ws2.Cells.Item(132, 9).Value = 4660
$ws2.Cells.Item(132, 10).Value = 4999
$ws2.Cells.Item(132, 11).Value = 13980
$ws2.Cells.Item(132, 12).Value = 14997
$ws2.Cells.Item(132, 13).Value = -11450
$ws2.Cells.Item(132, 14).Value = -20057

# ARM row 135
$ws2.Cells.Item(135, 8).Value = 100000
$ws2.Cells.Item(135, 9).Value = 0
$ws2.Cells.Item(135, 10).Value = 100000
$ws2.Cells.Item(135, 11).Value = 0
$ws2.Cells.Item(135, 12).Value = 100000
$ws2.Cells.Item(135, 14).Value = -110140

# ARM row 136
$ws2.Cells.Item(136, 8).Value = 3047.6667
$ws2.Cells.Item(136, 9).Value = 3047.6667
$ws2.Cells.Item(136, 10).Value = 0
$ws2.Cells.Item(136, 11).Value = 9143.000100000001
$ws2.Cells.Item(136, 12).Value = 0
$ws2.Cells.Item(136, 13).Value = -6593.000100000001

# BSM row 86
$ws3.Cells.Item(86, 8).Value = 7469.9
$ws3.Cells.Item(86, 9).Value = 2639.8
$ws3.Cells.Item(86, 10).Value = 12300
$ws3.Cells.Item(86, 11).Value = 2639.8
$ws3.Cells.Item(86, 12).Value = 12300
$ws3.Cells.Item(86, 13).Value = -1516.8
$ws3.Cells.Item(86, 14).Value = -14546

# BSM row 89
$ws3.Cells.Item(89, 8).Value = 7469.9
$ws3.Cells.Item(89, 9).Value = 2639.8
$ws3.Cells.Item(89, 10).Value = 12300
$ws3.Cells.Item(89, 11).Value = 13199
$ws3.Cells.Item(89, 12).Value = 61500
$ws3.Cells.Item(89, 13).Value = -7583
$ws3.Cells.Item(89, 14).Value = -72732

# BSM row 105
$ws3.Cells.Item(105, 8).Value = 1950.5
$ws3.Cells.Item(105, 9).Value = 1950.5
$ws3.Cells.Item(105, 10).Value = 0
$ws3.Cells.Item(105, 11).Value = 1950.5
$ws3.Cells.Item(105, 12).Value = 0
$ws3.Cells.Item(105, 13).Value = -203.5

# BSM row 107
$ws3.Cells.Item(107, 8).Value = 1835.5454
$ws3.Cells.Item(107, 9).Value = 1741.7142
$ws3.Cells.Item(107, 10).Value = 1999.75
$ws3.Cells.Item(107, 11).Value = 1741.7142
$ws3.Cells.Item(107, 12).Value = 1999.75
$ws3.Cells.Item(107, 13).Value = 178.2858000000001
$ws3.Cells.Item(107, 14).Value = -5839.75

# CRP row 22
$ws4.Cells.Item(22, 8).Value = 716.55554
$ws4.Cells.Item(22, 9).Value = 638.8570999999999
$ws4.Cells.Item(22, 10).Value = 988.5
$ws4.Cells.Item(22, 11).Value = 638.8570999999999
$ws4.Cells.Item(22, 12).Value = 988.5
$ws4.Cells.Item(22, 13).Value = -288.8570999999999
$ws4.Cells.Item(22, 14).Value = -1688.5

# CRP row 31
$ws4.Cells.Item(31, 8).Value = 2137.3125
$ws4.Cells.Item(31, 9).Value = 1585.5
$ws4.Cells.Item(31, 10).Value = 6000
$ws4.Cells.Item(31, 11).Value = 1585.5
$ws4.Cells.Item(31, 12).Value = 6000
$ws4.Cells.Item(31, 13).Value = -1290.5
$ws4.Cells.Item(31, 14).Value = -6590

# CRP row 34
$ws4.Cells.Item(34, 8).Value = 2137.3125
$ws4.Cells.Item(34, 9).Value = 1585.5
$ws4.Cells.Item(34, 10).Value = 6000
$ws4.Cells.Item(34, 11).Value = 1585.5
$ws4.Cells.Item(34, 12).Value = 6000
$ws4.Cells.Item(34, 13).Value = -1383.5
$ws4.Cells.Item(34, 14).Value = -6404

# CRP row 132
$ws4.Cells.Item(132, 8).Value = 3363
$ws4.Cells.Item(132, 9).Value = 2817.6667
$ws4.Cells.Item(132, 10).Value = 4999
$ws4.Cells.Item(132, 11).Value = 8453.000100000001
$ws4.Cells.Item(132, 12).Value = 14997
$ws4.Cells.Item(132, 13).Value = -5923.000100000001
$ws4.Cells.Item(132, 14).Value = -20057

# CUL row 33
$ws5.Cells.Item(33, 8).Value = 42.166668
$ws5.Cells.Item(33, 9).Value = 51.11111
$ws5.Cells.Item(33, 10).Value = 15.333333
$ws5.Cells.Item(33, 11).Value = 306.66666
$ws5.Cells.Item(33, 12).Value = 91.99999800000001
$ws5.Cells.Item(33, 13).Value = -23.66665999999998
$ws5.Cells.Item(33, 14).Value = -657.999998

# CUL row 127
$ws5.Cells.Item(127, 8).Value = 1200
$ws5.Cells.Item(127, 9).Value = 0
$ws5.Cells.Item(127, 10).Value = 1200
$ws5.Cells.Item(127, 11).Value = 0
$ws5.Cells.Item(127, 12).Value = 3600
$ws5.Cells.Item(127, 14).Value = -13520

# CUL row 131
$ws5.Cells.Item(131, 8).Value = 3267.5715
$ws5.Cells.Item(131, 9).Value = 2999.5
$ws5.Cells.Item(131, 10).Value = 3288.1924
$ws5.Cells.Item(131, 11).Value = 8998.5
$ws5.Cells.Item(131, 12).Value = 9864.5772
$ws5.Cells.Item(131, 13).Value = -3958.5
$ws5.Cells.Item(131, 14).Value = -19944.5772

# CUL row 134
$ws5.Cells.Item(134, 8).Value = 1578.75
$ws5.Cells.Item(134, 9).Value = 1578.75
$ws5.Cells.Item(134, 10).Value = 0
$ws5.Cells.Item(134, 11).Value = 4736.25
$ws5.Cells.Item(134, 12).Value = 0
$ws5.Cells.Item(134, 13).Value = 333.75

# GSM row 80
$ws6.Cells.Item(80, 8).Value = 3500
$ws6.Cells.Item(80, 9).Value = 3500
$ws6.Cells.Item(80, 10).Value = 0
$ws6.Cells.Item(80, 11).Value = 3500
$ws6.Cells.Item(80, 12).Value = 0
$ws6.Cells.Item(80, 13).Value = -2502
$ws6.Cells.Item(80, 14).ClearContents()

# GSM row 83
$ws6.Cells.Item(83, 8).Value = 3500
$ws6.Cells.Item(83, 9).Value = 3500
$ws6.Cells.Item(83, 10).Value = 0
$ws6.Cells.Item(83, 11).Value = 17500
$ws6.Cells.Item(83, 12).Value = 0
$ws6.Cells.Item(83, 13).Value = -12508
$ws6.Cells.Item(83, 14).ClearContents()

# GSM row 86
$ws6.Cells.Item(86, 8).Value = 0
$ws6.Cells.Item(86, 9).Value = 0
$ws6.Cells.Item(86, 10).Value = 0
$ws6.Cells.Item(86, 11).Value = 0
$ws6.Cells.Item(86, 12).Value = 0
$ws6.Cells.Item(86, 14).ClearContents()

# GSM row 89
$ws6.Cells.Item(89, 8).Value = 0
$ws6.Cells.Item(89, 9).Value = 0
$ws6.Cells.Item(89, 10).Value = 0
$ws6.Cells.Item(89, 11).Value = 0
$ws6.Cells.Item(89, 12).Value = 0
$ws6.Cells.Item(89, 14).ClearContents()

# GSM row 126
$ws6.Cells.Item(126, 8).Value = 5249.5
$ws6.Cells.Item(126, 9).Value = 5399.4
$ws6.Cells.Item(126, 10).Value = 4500
$ws6.Cells.Item(126, 11).Value = 16198.2
$ws6.Cells.Item(126, 12).Value = 13500
$ws6.Cells.Item(126, 13).Value = -13728.2
$ws6.Cells.Item(126, 14).Value = -18440

# WVR row 82
$ws8.Cells.Item(82, 8).Value = 0
$ws8.Cells.Item(82, 9).Value = 0
$ws8.Cells.Item(82, 10).Value = 0
$ws8.Cells.Item(82, 11).Value = 0
$ws8.Cells.Item(82, 12).Value = 0
$ws8.Cells.Item(82, 14).ClearContents()

# WVR row 85
$ws8.Cells.Item(85, 8).Value = 0
$ws8.Cells.Item(85, 9).Value = 0
$ws8.Cells.Item(85, 10).Value = 0
$ws8.Cells.Item(85, 11).Value = 0
$ws8.Cells.Item(85, 12).Value = 0
$ws8.Cells.Item(85, 14).ClearContents()

# WVR row 122
$ws8.Cells.Item(122, 8).Value = 1138.7778
$ws8.Cells.Item(122, 9).Value = 1138.7778
$ws8.Cells.Item(122, 10).Value = 0
$ws8.Cells.Item(122, 11).Value = 3416.3334
$ws8.Cells.Item(122, 12).Value = 0
$ws8.Cells.Item(122, 13).Value = -966.3334000000004
$ws8.Cells.Item(122, 14).ClearContents()

# WVR row 132
$ws8.Cells.Item(132, 8).Value = 4185.923
$ws8.Cells.Item(132, 9).Value = 3841.111
$ws8.Cells.Item(132, 10).Value = 4961.75
$ws8.Cells.Item(132, 11).Value = 11523.333
$ws8.Cells.Item(132, 12).Value = 14885.25
$ws8.Cells.Item(132, 13).Value = -8993.332999999999
$ws8.Cells.Item(132, 14).Value = -19945.25
